# Populate "output.xlsx" (Sheet1) with the QR-code reader results:
#   - adds a new "qr_code_text" column (K) with the decoded QR payload for
#     every slip, matching the header style already used by A1:J1
#   - fills in the previously-unread rows (file/sender/receiver/etc.) and
#     appends two additional parsed slips (Mymoo.txt, SCB.txt), growing the
#     sheet from 11 to 13 data+header rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("file", "bank", "sender_name", "sender_bank", "sender_acc", "receiver_name", "receiver_bank", "receiver_acc", "amount", "date", "qr_code_text"),
  @("BKB.txt", "bangkok", "น.ส. ฺทัดธนัณ", "กรุงเทพ", "266-4 xxx864", "น.ส. ภัทราวดี  จันทร์ปากน้ำ", "พร้อมเพย์", "266-4 xxx864", "20.00", "16 มิ.ย. 68, 19:26", "004600060000010103002022520250616192655230052935085102TH91045015"),
  @("BKB2000.txt", "bangkok", "น.ส. ทัดธนัณ", "กรุงเทพ", "266-4 xxx864", "น.ส. ทัดธนัณ ขุนแดง", "ธนาคารกรุงไทย", "266-4 xxx864", "2,000.00", "06 เม.ย. 68, 12:32", "004600060000010103002022520250406123215230058366085102TH91040A46"),
  @("d25b87bd-1da5-4df7-b115-f0b8b5596d1a.txt", "gsb", "นาย อดิเทพ แบบเหมือน", "ธนาคารออมสิน", "0204xxxx1195", ". นายสัญญา ศรีวรรณ.", "พร้อมเพย์", "3909xxxxx4703", "25.00", "15ก.ค. 2568 17:32", "0045000600000101030300224519617039587I000007B97905102TH9104D7A2"),
  @("K+14000.txt", "kbank", "น.ส. ซิลมีย์ ป.", "ธนาคารกสิกรไทย", "xxx-x-x2730-x", "น.ส. ภักดิพระสามารถ ธ.กสิกรไทย", "ธ.กสิกรไทย", "xxx-x-x5833-x", "14,000.00", "19 เม.ย. 68 18:07 น.", "0041000600000101030040220015109180730BTF085585102TH91042F11"),
  @("K+2000.txt", "kbank", "น.ส. ซิลมีย์ ป.", "ธนาคารกสิกรไทย", "xxx-x-x2730-x", "น.ส. ซิลมีย์ ปะหนัน", "ธ.กสิกรไทย", "xxx-x-x6833-x", "2,000.00", "30 พ.ค. 68 06:03 น.", "0041000600000101030040220015150060303ATF010635102TH91047FF0"),
  @("K+kitty.txt", "kbank", "น.ส. ภักดิพรส ธ.กสิกรไทย", "ธนาคารกสิกรไทย", "x-x5833-x", "น.ส. ภักดิพร สามารถ", "ธ.กสิกรไทย", "xxx-x-x3478-x", "0.84", "16 ก.ค. 68  23:14 น.", "0041000600000101030040220015197231427CTF016705102TH9104D930"),
  @("K+snoopy.txt", "kbank", "ไม่พบชื่อผู้โอน", "ธนาคารกสิกรไทย", "xxx-x-x5833-x", "ไม่พบชื่อผู้รับ", "ไม่พบชื่อธนาคารผู้รับ", "ไม่พบเลขบัญชีผู้รับ", "3,325.00", "ไม่พบวันที่", "0041000600000101030040220015195192032APM136485102TH91041E3C"),
  @("KB.txt", "kbank", "น.ส. ซิลมีย์ ป.", "ธนาคารกสิกรไทย", "xxx-x-x2730-x", "น.ส. ภักดิพระสามารถ ธ.กสิกรไทย", "ธ.กสิกรไทย", "xxx-x-x5833-x", "14,000.00", "19 เม.ย. 68 18:07 น.", "0041000600000101030040220015109180730BTF085585102TH91042F11"),
  @("KTB.txt", "krungthai", "ซิลมีย์", "กรุงไทย", "x-xx882-2", "น.ส. ซิลมีย์ ปะหนัน", "กสิกรไทย", "x-xx833-8", "500.00", "20 ก.พ. 2568 15:19", "0038000600000101030060217A78990c2774ba4b195102TH91045737"),
  @("Mymoo.txt", "gsb", "นาย อดิเทพ แบบเหมือน", "ธนาคารออมสิน", "0204xxxx1195", "นาย อัฟเฟนดีย์สาอิ", "เติมเงินพร้อมเพย์", "4999xxxxx0649", "30.00", "15ก.ค. 2568 17:59", "0045000600000101030300224519617054995I000007B97905102TH910450FA"),
  @("SCB.txt", "scb", "นางสาว ซิลมีย์ ปะหนัน", "ไทยพาณิชย์", "xxx xxx189-6", "น.ส.ศิริวรรณ หมะหวังเอียด", "ไม่ระบุ", "xxx xxx-9283", "40.00", "25 มิ.ย. 2568 17:53", "0046000600000101030140225202506258EWX6Ejo07FyhP4ev5102TH91048888"),
  @("SCB3500.txt", "scb", "นางสาว ซิลมีย์ ปะหนัน", "ไทยพาณิชย์", "x482", "น.ส. ซิลมีย์ ปะหนัน", "ไม่ระบุ", "x-8338", "3,500.00", "01 ก.ค. 2568 10:15", "0046000600000101030140225202507018x482jp7EtCLE2oY15102TH9104F5B4")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $rowNum = $i + 1
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $colNum = $j + 1
    $cell = $ws.Cells.Item($rowNum, $colNum)
    $text = $row[$j]
    if ($colNum -eq 9 -and $rowNum -gt 1) {
      # "amount" column: force text so numeric-looking values like "20.00" or
      # "2,000.00" are kept verbatim instead of being parsed into numbers, then
      # clear the formatting again so the cell keeps its original (no) style.
      $cell.NumberFormat = "@"
      $cell.Value = $text
      $cell.ClearFormats()
    } else {
      $cell.Value = $text
    }
  }
}

# Give the new K1 header cell ("qr_code_text") the same bold/bordered/centered
# style already used by the rest of row 1 (copy format from J1 "date").
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "done"
